# pico-printable.pptx edit:
#  - slide 1: remove the 3 "bottom-row"/odd pictures, shift the 3 remaining
#    pictures to sit (almost) flush with the top of the slide, then duplicate
#    each of them to rebuild a second (bottom) row of pictures.
#  - add speaker notes infrastructure (notes placeholder) for slide 1.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# EMU -> Points helper (Shape.Left/Top/Width/Height are expressed in points).
$EMU_PER_POINT = 12700

function EmuToPt($emu) {
    return [double]$emu / $EMU_PER_POINT
}

# --- 1. Delete the three pictures that are being removed -------------------
# (identified by their original shape Id in the source deck)
$idsToDelete = @(5, 8, 10)
for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
    $shp = $s.Shapes.Item($i)
    if ($idsToDelete -contains $shp.Id) {
        $shp.Delete()
    }
}

# --- 2. Reposition the three pictures that remain ---------------------------
# Map: original shape Id -> new (Left, Top) in EMU
$newPos = @{
    7  = @(324944, 1)
    9  = @(3240356, 1)
    11 = @(6155768, 11432)
}

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($newPos.ContainsKey($shp.Id)) {
        $pos = $newPos[$shp.Id]
        $shp.Left = EmuToPt $pos[0]
        $shp.Top  = EmuToPt $pos[1]
    }
}

# --- 3. Duplicate the three remaining pictures to build a second row -------
# Map: original shape Id -> (new duplicate's name, new Left/Top in EMU)
$dupInfo = @{
    7  = @{ Name = "Picture 1"; Pos = @(324944, 3389377) }
    9  = @{ Name = "Picture 2"; Pos = @(3240356, 3389377) }
    11 = @{ Name = "Picture 3"; Pos = @(6155768, 3400808) }
}

# Snapshot current shapes first (Duplicate adds to the end of the collection,
# so iterating live would also re-visit the newly added duplicates).
$sourceIds = @()
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sourceIds += $s.Shapes.Item($i).Id
}

foreach ($srcId in $sourceIds) {
    if ($dupInfo.ContainsKey($srcId)) {
        # locate the source shape again by Id
        for ($i = 1; $i -le $s.Shapes.Count; $i++) {
            $cand = $s.Shapes.Item($i)
            if ($cand.Id -eq $srcId) {
                $info = $dupInfo[$srcId]
                $dup = $cand.Duplicate()
                $dup.Name = $info.Name
                $dup.Left = EmuToPt $info.Pos[0]
                $dup.Top  = EmuToPt $info.Pos[1]
                break
            }
        }
    }
}

# --- 4. Speaker notes: create the notes placeholder for slide 1 ------------
$notes = $s.NotesPage
$notesBody = $notes.Shapes.AddPlaceholder(2)
